# Apply BOM fixes:
# 1. Rename sheet from "BOM_ESP32-Wiegand-Hat_ESP32-Wie" to "BOM"
# 2. Fix Quantity for row 2 (PZ254V-11-01P) from 2 to 1
# 3. Fix Designator for row 2 from "5V,12V" to "12V"
# 4. Fix Footprint casing for row 7 (Level Shifter) from "level shifter footprint" to "Level Shifter Footprint"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "BOM"

# B2 needs to become "1" (text), matching the text-typed "1" already present
# in B3. Copy/PasteSpecial (values) from B3 keeps it stored as text instead
# of Excel auto-converting a plain "1" into a number.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D2").Value = "12V"
$ws.Range("E7").Value = "Level Shifter Footprint"
